$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J15").Value = 17
$ws.Range("K15").Value = -70.588235294117
$ws.Range("M15").Value = -82.758620689655
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 77.777777777777
$ws.Range("I16").Value = 120
$ws.Range("J16").Value = 119
$ws.Range("K16").Value = 0.840336134453
$ws.Range("L16").Value = 26.315789473684
$ws.Range("M16").Value = -46.188340807174
$ws.Range("N16").Value = -81.220657276995
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 210
$ws.Range("J17").Value = 202
$ws.Range("K17").Value = 3.960396039603
$ws.Range("L17").Value = 17.977528089887
$ws.Range("M17").Value = 30.434782608695
$ws.Range("N17").Value = -66.453674121405
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 57.142857142857
$ws.Range("I18").Value = 67
$ws.Range("J18").Value = 66
$ws.Range("K18").Value = 1.515151515151
$ws.Range("L18").Value = -34.313725490196
$ws.Range("M18").Value = -32.323232323232
$ws.Range("N18").Value = -90.482954545454
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 19.047619047619
$ws.Range("I19").Value = 316
$ws.Range("J19").Value = 318
$ws.Range("K19").Value = -0.628930817610
$ws.Range("L19").Value = 7.849829351535
$ws.Range("M19").Value = 137.593984962406
$ws.Range("N19").Value = -1.557632398753
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -22.222222222222
$ws.Range("I20").Value = 104
$ws.Range("J20").Value = 77
$ws.Range("K20").Value = 35.064935064935
$ws.Range("L20").Value = 85.714285714285
$ws.Range("M20").Value = 131.111111111111
$ws.Range("N20").Value = -64.625850340136
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 58.333333333333
$ws.Range("F21").Value = 69
$ws.Range("G21").Value = 57
$ws.Range("H21").Value = 21.052631578947
$ws.Range("I21").Value = 825
$ws.Range("J21").Value = 803
$ws.Range("K21").Value = 2.739726027397
$ws.Range("L21").Value = 12.551159618008
$ws.Range("M21").Value = 18.534482758620
$ws.Range("N21").Value = -69.216417910447
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("M22").Value = 115.384615384615
$ws.Range("M23").Value = 116.666666666667
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = 37.5
$ws.Range("F24").Value = 47
$ws.Range("G24").Value = 51
$ws.Range("H24").Value = -7.843137254901
$ws.Range("I24").Value = 703
$ws.Range("J24").Value = 1051
$ws.Range("K24").Value = -33.111322549952
$ws.Range("L24").Value = -10.786802030456
$ws.Range("M24").Value = 97.471910112359
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 400
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = 88.888888888888
$ws.Range("I25").Value = 301
$ws.Range("J25").Value = 315
$ws.Range("K25").Value = -4.444444444444
$ws.Range("L25").Value = -0.660066006600
$ws.Range("M25").Value = -35.407725321888
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 24
$ws.Range("K26").Value = -41.666666666666
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 31
$ws.Range("J27").Value = 51
$ws.Range("K27").Value = -39.215686274509
$ws.Range("L27").Value = 6.896551724137
$ws.Range("M28").Value = -61.538461538461
$ws.Range("M29").Value = -54.545454545454

# Cells that change from blank/text dash to numeric (style 14 -> 15/16)
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value = -100
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("G15").Value = 1
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H15").Value = -100

$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 1
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E26").Value = -100

$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 2
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E27").Value = -50

# C22: numeric (style 15) -> blank dash text (style 14, shared string "0")
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("A22").Copy()
$ws.Range("C22").PasteSpecial(-4122)

# Rich-text shared strings: in-place run edits via Characters
$ws.Range("A8").Characters(21,2).Text = "47"
$ws.Range("C9").Characters(27,10).Text = "11/20/2023"
$ws.Range("C9").Characters(48,10).Text = "11/26/2023"
